# Handback status report regeneration:
#  - the in-flight handback file 4abcd790-...  was renamed/superseded by
#    6e6e9a67-c2bc-4596-a7dc-042ca4ece98f (new hashes/timestamps), and a
#    second handback file e5e82743-bf71-404a-8a02-c0f5851885ee was added
#    as a brand-new row on every sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")

# Row 2: 4abcd790... -> 6e6e9a67... (same slot, new file id + new date)
$ov.Range("A2").Value = "6e6e9a67-c2bc-4596-a7dc-042ca4ece98f.md"
$ov.Range("B2").Hyperlinks.Delete()
$ov.Hyperlinks.Add($ov.Range("B2"), `
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/80132af0e5a00bcc1821f91c5a1a4c202a4f6821/e2e/6e6e9a67-c2bc-4596-a7dc-042ca4ece98f.md", `
    "", "", "e2e\6e6e9a67-c2bc-4596-a7dc-042ca4ece98f.md") | Out-Null
$ov.Range("G2").Value = "2016-08-25 06:59:58"
$ov.Range("G2").NumberFormat = "yyyy-mm-dd HH:mm:ss"

# Row 3: brand-new handback file e5e82743...
$ovTable = $ov.ListObjects.Item(1)
$ovTable.ListRows.Add() | Out-Null
$ov.Range("A3").Value = "e5e82743-bf71-404a-8a02-c0f5851885ee.md"
$ov.Range("C3").Value = ".md"
$ov.Range("E3").Value = "Handed back: in sync with en-US"
$ov.Range("F3").Value = "Handed back: in sync with en-US"
$ov.Range("G3").Value = "2016-08-25 06:59:58"
$ov.Range("G3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ov.Hyperlinks.Add($ov.Range("B3"), `
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/80132af0e5a00bcc1821f91c5a1a4c202a4f6821/e2e/e5e82743-bf71-404a-8a02-c0f5851885ee.md", `
    "", "", "e2e\e5e82743-bf71-404a-8a02-c0f5851885ee.md") | Out-Null

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

# Row 2: 4abcd790... -> 6e6e9a67... (new hash + refreshed timestamps)
$zh.Range("A2").Value = "6e6e9a67-c2bc-4596-a7dc-042ca4ece98f.md"
$zh.Range("A2").Hyperlinks.Delete()
$zh.Hyperlinks.Add($zh.Range("A2"), `
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/80132af0e5a00bcc1821f91c5a1a4c202a4f6821/e2e/6e6e9a67-c2bc-4596-a7dc-042ca4ece98f.md", `
    "", "", "6e6e9a67-c2bc-4596-a7dc-042ca4ece98f.md") | Out-Null
$zh.Range("G2").Value = "6e6e9a67-c2bc-4596-a7dc-042ca4ece98f.27c1ee65b98a7ee1dd76788151f3087b9734550b.zh-cn.xlf"
$zh.Range("H2").Value = "2016-08-25 06:59:52"
$zh.Range("I2").Value = "6e6e9a67-c2bc-4596-a7dc-042ca4ece98f.md"
$zh.Range("I2").Hyperlinks.Delete()
$zh.Hyperlinks.Add($zh.Range("I2"), `
    "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/776ea6a1ef52a990c317ec83d919b16b0bbc0502/e2e/6e6e9a67-c2bc-4596-a7dc-042ca4ece98f.md", `
    "", "", "6e6e9a67-c2bc-4596-a7dc-042ca4ece98f.md") | Out-Null
$zh.Range("J2").Value = "6e6e9a67-c2bc-4596-a7dc-042ca4ece98f.27c1ee65b98a7ee1dd76788151f3087b9734550b.zh-cn.xlf"
$zh.Range("K2").Value = "2016-08-25 07:00:35"

# Row 3: brand-new handback file e5e82743...
$zhTable = $zh.ListObjects.Item(1)
$zhTable.ListRows.Add() | Out-Null
$zh.Range("A3").Value = "e5e82743-bf71-404a-8a02-c0f5851885ee.md"
$zh.Hyperlinks.Add($zh.Range("A3"), `
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/80132af0e5a00bcc1821f91c5a1a4c202a4f6821/e2e/e5e82743-bf71-404a-8a02-c0f5851885ee.md", `
    "", "", "e5e82743-bf71-404a-8a02-c0f5851885ee.md") | Out-Null
$zh.Range("B3").Value = ".md"
$zh.Range("C3").Value = "Handed back: in sync with en-US"
$zh.Range("D3").Value = "e2e"
$zh.Range("E3").Value = "ht"
$zh.Range("F3").Value = "True"
$zh.Range("G3").Value = "e5e82743-bf71-404a-8a02-c0f5851885ee.2fa26be12d8952b6ce93288ba391e56e600b07c7.zh-cn.xlf"
$zh.Range("H3").Value = "2016-08-25 06:59:52"
$zh.Range("I3").Value = "e5e82743-bf71-404a-8a02-c0f5851885ee.md"
$zh.Hyperlinks.Add($zh.Range("I3"), `
    "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/776ea6a1ef52a990c317ec83d919b16b0bbc0502/e2e/e5e82743-bf71-404a-8a02-c0f5851885ee.md", `
    "", "", "e5e82743-bf71-404a-8a02-c0f5851885ee.md") | Out-Null
$zh.Range("J3").Value = "e5e82743-bf71-404a-8a02-c0f5851885ee.2fa26be12d8952b6ce93288ba391e56e600b07c7.zh-cn.xlf"
$zh.Range("K3").Value = "2016-08-25 07:00:35"
$zh.Range("M3").Value = "True"
$zh.Range("O3").Value = "False"

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

# Row 2: 4abcd790... -> 6e6e9a67... (new hash + refreshed timestamps)
$de.Range("A2").Value = "6e6e9a67-c2bc-4596-a7dc-042ca4ece98f.md"
$de.Range("A2").Hyperlinks.Delete()
$de.Hyperlinks.Add($de.Range("A2"), `
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/80132af0e5a00bcc1821f91c5a1a4c202a4f6821/e2e/6e6e9a67-c2bc-4596-a7dc-042ca4ece98f.md", `
    "", "", "6e6e9a67-c2bc-4596-a7dc-042ca4ece98f.md") | Out-Null
$de.Range("G2").Value = "6e6e9a67-c2bc-4596-a7dc-042ca4ece98f.27c1ee65b98a7ee1dd76788151f3087b9734550b.de-de.xlf"
$de.Range("I2").Value = "6e6e9a67-c2bc-4596-a7dc-042ca4ece98f.md"
$de.Range("I2").Hyperlinks.Delete()
$de.Hyperlinks.Add($de.Range("I2"), `
    "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/482b5487603da6abb051c9a532902db484d56dfa/e2e/6e6e9a67-c2bc-4596-a7dc-042ca4ece98f.md", `
    "", "", "6e6e9a67-c2bc-4596-a7dc-042ca4ece98f.md") | Out-Null
$de.Range("J2").Value = "6e6e9a67-c2bc-4596-a7dc-042ca4ece98f.27c1ee65b98a7ee1dd76788151f3087b9734550b.de-de.xlf"
$de.Range("K2").Value = "2016-08-25 07:00:43"

# Row 3: brand-new handback file e5e82743...
$deTable = $de.ListObjects.Item(1)
$deTable.ListRows.Add() | Out-Null
$de.Range("A3").Value = "e5e82743-bf71-404a-8a02-c0f5851885ee.md"
$de.Hyperlinks.Add($de.Range("A3"), `
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/80132af0e5a00bcc1821f91c5a1a4c202a4f6821/e2e/e5e82743-bf71-404a-8a02-c0f5851885ee.md", `
    "", "", "e5e82743-bf71-404a-8a02-c0f5851885ee.md") | Out-Null
$de.Range("B3").Value = ".md"
$de.Range("C3").Value = "Handed back: in sync with en-US"
$de.Range("D3").Value = "e2e"
$de.Range("E3").Value = "ht"
$de.Range("F3").Value = "True"
$de.Range("G3").Value = "e5e82743-bf71-404a-8a02-c0f5851885ee.2fa26be12d8952b6ce93288ba391e56e600b07c7.de-de.xlf"
$de.Range("H3").Value = "2016-08-25 06:59:58"
$de.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$de.Range("I3").Value = "e5e82743-bf71-404a-8a02-c0f5851885ee.md"
$de.Hyperlinks.Add($de.Range("I3"), `
    "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/482b5487603da6abb051c9a532902db484d56dfa/e2e/e5e82743-bf71-404a-8a02-c0f5851885ee.md", `
    "", "", "e5e82743-bf71-404a-8a02-c0f5851885ee.md") | Out-Null
$de.Range("J3").Value = "e5e82743-bf71-404a-8a02-c0f5851885ee.2fa26be12d8952b6ce93288ba391e56e600b07c7.de-de.xlf"
$de.Range("K3").Value = "2016-08-25 07:00:43"
$de.Range("M3").Value = "True"
$de.Range("O3").Value = "False"

Write-Output "Handback status report regenerated."
